# Apply crypto price/volume updates per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.192.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.33%  "

$ws.Range("D3").Value = "'2.425.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'308.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.62%  "

$ws.Range("D6").Value = "'100.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.27%  "

$ws.Range("D7").Value = "'0.512"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.500"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").Value = "'35.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("D11").Value = "'0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("E12").Value = "  +1.97%  "

$ws.Range("D13").Value = "'18.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("D15").Value = "'2.802.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "

$ws.Range("D16").Value = "'2.415.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "'0.834"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.18%  "

$ws.Range("D18").Value = "'44.138.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("E19").Value = "  +1.21%  "

$ws.Range("D20").Value = "'6.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("D21").Value = "'0.0₃0904"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.67%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'240.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.94%  "

$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.48%  "

$ws.Range("E25").Value = "  +1.05%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").Value = "'25.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "

$ws.Range("D28").Value = "'2.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.39%  "

$ws.Range("D29").Value = "'9.60"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").Value = "'32.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.84%  "

$ws.Range("D31").Value = "'18.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.40%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.48%  "

$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "'0.116"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.37%  "

$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("E35").Value = "  +2.37%  "

$ws.Range("E36").Value = "  +3.11%  "

$ws.Range("D37").Value = "'4.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.82%  "

$ws.Range("D38").Value = "'129.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +22.11%  "

$ws.Range("D39").Value = "'2.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "

$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").Value = "'21.58"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.38%  "

$ws.Range("E43").Value = "  +2.44%  "

$ws.Range("D44").Value = "'1.955.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("E46").Value = "  +4.12%  "

$ws.Range("D47").Value = "'9.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.15%  "

$ws.Range("D48").Value = "'1.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.27%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "'2.667.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.42%  "

$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "'53.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'73.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "

Write-Output "Updated crypto rows."